# Applies the daily "remaining days" decrement update.
# For every data row (2..99) in column E (剩余/remaining), subtract 1.
# If the remaining count would reach 0 (i.e. it was 1), instead reset it
# to 10 and push the start date (column F) forward by 10 days, starting
# a new cycle.
# Row 36 has a corrupted/out-of-range date value and is left untouched,
# matching the source data which also skips it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        continue
    }

    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F

    $eValue = $eCell.Value2()

    if ($null -eq $eValue) {
        continue
    }

    $eValue = [int]$eValue

    if ($eValue -eq 1) {
        $eCell.Value = 10
        $fValue = [int]$fCell.Value2()
        $fCell.Value = $fValue + 10
    } else {
        $eCell.Value = $eValue - 1
    }
}
